$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6321.3335
$ws.Range("I18").Value = 3318.4167
$ws.Range("J18").Value = 18333
$ws.Range("K18").Value = 3318.4167
$ws.Range("L18").Value = 18333
$ws.Range("M18").Value = -3034.4167
$ws.Range("N18").Value = -18901
$ws.Range("H53").Value = 315.1
$ws.Range("I53").Value = 129.5
$ws.Range("J53").Value = 361.5
$ws.Range("K53").Value = 129.5
$ws.Range("L53").Value = 361.5
$ws.Range("M53").Value = 507.5
$ws.Range("N53").Value = -1635.5
$ws.Range("H100").Value = 2312
$ws.Range("I100").Value = 2431.2727
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 2431.2727
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -1890.2727
$ws.Range("N100").Value = -2082
$ws.Range("H137").Value = 24181.227
$ws.Range("I137").Value = 46383
$ws.Range("J137").Value = 17705.709
$ws.Range("K137").Value = 139149
$ws.Range("L137").Value = 53117.12699999999
$ws.Range("M137").Value = -136599
$ws.Range("N137").Value = -58217.12699999999
$ws.Range("H138").Value = 23452.646
$ws.Range("I138").Value = 2026.409
$ws.Range("K138").Value = 6079.227000000001
$ws.Range("M138").Value = -939.2270000000008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36395.13
$ws.Range("I32").Value = 39711.32
$ws.Range("K32").Value = 39711.32
$ws.Range("M32").Value = -39424.32
$ws.Range("H34").Value = 144548
$ws.Range("J34").Value = 109002.8
$ws.Range("L34").Value = 109002.8
$ws.Range("N34").Value = -109544.8
$ws.Range("H61").Value = 8847.474
$ws.Range("I61").Value = 937.61536
$ws.Range("K61").Value = 937.61536
$ws.Range("M61").Value = -725.61536
$ws.Range("H74").Value = 360127.34
$ws.Range("I74").Value = 429482.22
$ws.Range("K74").Value = 429482.22
$ws.Range("M74").Value = -428608.22
$ws.Range("H77").Value = 360127.34
$ws.Range("I77").Value = 429482.22
$ws.Range("K77").Value = 2147411.1
$ws.Range("M77").Value = -2143043.1
$ws.Range("H132").Value = 2956
$ws.Range("J132").Value = 4658
$ws.Range("L132").Value = 13974
$ws.Range("N132").Value = -19034
$ws.Range("H136").Value = 8847.474
$ws.Range("I136").Value = 937.61536
$ws.Range("K136").Value = 2812.84608
$ws.Range("M136").Value = -262.8460800000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 10249
$ws.Range("I36").Value = 2333
$ws.Range("J36").Value = 33997
$ws.Range("K36").Value = 2333
$ws.Range("L36").Value = 33997
$ws.Range("M36").Value = -1799
$ws.Range("N36").Value = -35065
$ws.Range("H105").Value = 1929.7241
$ws.Range("I105").Value = 1620.8077
$ws.Range("K105").Value = 1620.8077
$ws.Range("M105").Value = 126.1922999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2555.7144
$ws.Range("I16").Value = 2663.3333
$ws.Range("J16").Value = 2475
$ws.Range("K16").Value = 2663.3333
$ws.Range("L16").Value = 2475
$ws.Range("M16").Value = -2376.3333
$ws.Range("N16").Value = -3049
$ws.Range("H94").Value = 1720.2858
$ws.Range("I94").Value = 1194.909
$ws.Range("J94").Value = 2298.2
$ws.Range("K94").Value = 1194.909
$ws.Range("L94").Value = 2298.2
$ws.Range("M94").Value = -743.9090000000001
$ws.Range("N94").Value = -3200.2
$ws.Range("H99").Value = 1739.5
$ws.Range("I99").Value = 1320.6666
$ws.Range("J99").Value = 2158.3333
$ws.Range("K99").Value = 1320.6666
$ws.Range("L99").Value = 2158.3333
$ws.Range("M99").Value = 177.3334
$ws.Range("N99").Value = -5154.3333
$ws.Range("H113").Value = 2555.7144
$ws.Range("I113").Value = 2663.3333
$ws.Range("J113").Value = 2475
$ws.Range("K113").Value = 2663.3333
$ws.Range("L113").Value = 2475
$ws.Range("M113").Value = -493.3332999999998
$ws.Range("N113").Value = -6815
$ws.Range("H126").Value = 1739.5
$ws.Range("I126").Value = 1320.6666
$ws.Range("J126").Value = 2158.3333
$ws.Range("K126").Value = 3961.9998
$ws.Range("L126").Value = 6474.999899999999
$ws.Range("M126").Value = -1491.9998
$ws.Range("N126").Value = -11414.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 701.6667
$ws.Range("I33").Value = 477.5
$ws.Range("K33").Value = 2865
$ws.Range("M33").Value = -2582
$ws.Range("H41").Value = 2266.6667
$ws.Range("I41").Value = 900
$ws.Range("K41").Value = 2700
$ws.Range("M41").Value = -2362
$ws.Range("H63").Value = 6863.636
$ws.Range("H66").Value = 6863.636
$ws.Range("H87").Value = 11428.536
$ws.Range("I87").Value = 9999.666999999999
$ws.Range("J87").Value = 11600
$ws.Range("K87").Value = 29999.001
$ws.Range("L87").Value = 34800
$ws.Range("M87").Value = -28751.001
$ws.Range("N87").Value = -37296
$ws.Range("H90").Value = 11428.536
$ws.Range("I90").Value = 9999.666999999999
$ws.Range("J90").Value = 11600
$ws.Range("K90").Value = 89997.003
$ws.Range("L90").Value = 104400
$ws.Range("M90").Value = -83757.003
$ws.Range("N90").Value = -116880
$ws.Range("H137").Value = 4365.857
$ws.Range("I137").Value = 4112.2
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 12336.6
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -7236.599999999999
$ws.Range("N137").Value = -25200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 39090.91
$ws.Range("J123").Value = 39090.91
$ws.Range("L123").Value = 39090.91
$ws.Range("N123").Value = -43990.91
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2747
$ws.Range("I61").Value = 2829.3333
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 2829.3333
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -2627.3333
$ws.Range("N61").Value = -2904
$ws.Range("H68").Value = 3404.5625
$ws.Range("I68").Value = 2368.625
$ws.Range("J68").Value = 4440.5
$ws.Range("K68").Value = 2368.625
$ws.Range("L68").Value = 4440.5
$ws.Range("M68").Value = -1619.625
$ws.Range("N68").Value = -5938.5
$ws.Range("H71").Value = 3404.5625
$ws.Range("I71").Value = 2368.625
$ws.Range("J71").Value = 4440.5
$ws.Range("K71").Value = 11843.125
$ws.Range("L71").Value = 22202.5
$ws.Range("M71").Value = -8099.125
$ws.Range("N71").Value = -29690.5
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H93").Value = 2999.2222
$ws.Range("I93").Value = 2998.5
$ws.Range("J93").Value = 2999.8
$ws.Range("K93").Value = 2998.5
$ws.Range("L93").Value = 2999.8
$ws.Range("M93").Value = -1750.5
$ws.Range("N93").Value = -5495.8
$ws.Range("H100").Value = 2922.7778
$ws.Range("I100").Value = 2815
$ws.Range("K100").Value = 2815
$ws.Range("M100").Value = -2274
$ws.Range("H113").Value = 2747
$ws.Range("I113").Value = 2829.3333
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 2829.3333
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -659.3332999999998
$ws.Range("N113").Value = -6840
$ws.Range("H136").Value = 5970.385
$ws.Range("I136").Value = 5870.1113
$ws.Range("J136").Value = 6196
$ws.Range("K136").Value = 17610.3339
$ws.Range("L136").Value = 18588
$ws.Range("M136").Value = -15060.3339
$ws.Range("N136").Value = -23688

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("H113").Value = 1057.6296
$ws.Range("I113").Value = 1030.421
$ws.Range("J113").Value = 1122.25
$ws.Range("K113").Value = 3091.263
$ws.Range("L113").Value = 3366.75
$ws.Range("M113").Value = -921.2629999999999
$ws.Range("N113").Value = -7706.75
$ws.Range("H133").Value = 56250
$ws.Range("J133").Value = 56250
$ws.Range("L133").Value = 56250
$ws.Range("N133").Value = -66370
